$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update A4 with the new phone number
$ws.Range("A4").Value = 8273873833

# Remove rows 5, 6 and 7 entirely (shifting cells up / shrinking used range)
$ws.Range("A5:D7").Delete()

# Update the active selection to A5 (now the first empty row)
$ws.Range("A5").Select()
